# Update the "tasks" worksheet:
#  - mark "Master Data Test Case" (row 19) as Complete
#  - add new rows for Transaction Data Test Case, Core test data,
#    Master Data Test Data and Transaction Data Test Data (rows 20-23)
#  - grow the table / dimension accordingly
#  - update the selection to the last edited cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: Master Data Test Case status changes from "In process" to "Complete"
$ws.Range("C19").Value = "Complete"

# Duplicate row 19 four times (as rows 20-23), preserving all formatting/styles
$ws.Rows("19").Copy()
$ws.Rows("20").Insert(-4121)
$ws.Rows("19").Copy()
$ws.Rows("20").Insert(-4121)
$ws.Rows("19").Copy()
$ws.Rows("20").Insert(-4121)
$ws.Rows("19").Copy()
$ws.Rows("20").Insert(-4121)

# Row 20: Transaction Data Test Case
$ws.Range("A20").Value = "'        3.1.2"
$ws.Range("B20").Value = "Transaction Data Test Case"
$ws.Range("C20").Value = "In process"

# Row 21: Core test data
$ws.Range("A21").Value = "'    3.2"
$ws.Range("B21").Value = "Core test data"
$ws.Range("C21").Value = "In process"

# Row 22: Master Data Test Data
$ws.Range("A22").Value = "'        3.2.1"
$ws.Range("B22").Value = "Master Data Test Data"
$ws.Range("C22").Value = "In process"

# Row 23: Transaction Data Test Data
$ws.Range("A23").Value = "'        3.2.2"
$ws.Range("B23").Value = "Transaction Data Test Data"
$ws.Range("C23").Value = "In process"

# Grow the worksheet's table (ListObject) to cover the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F23"))

# Scroll / select to match the final view state
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C23").Select()
